# Update the cached "datetime8" field text on the Date placeholders of the
# Handout Master, Notes Master, and every Notes Page that has one, from the
# stale "9/17/18 3:44 PM" capture to the refreshed "3/4/19 8:29 PM" one.

$p = $ppt.ActivePresentation
$newDate = "3/4/19 8:29 PM"

# Handout Master - "Date Placeholder 6" is the 2nd shape.
$hm = $p.HandoutMaster
$hm.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# Notes Master - "Date Placeholder 10" is the 4th shape.
$nm = $p.NotesMaster
$nm.Shapes.Item(4).TextFrame.TextRange.Text = $newDate

# Every slide's Notes Page - "Date Placeholder 5" is the 5th shape, when present.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.HasNotesPage) {
        $notesPage = $slide.NotesPage
        $notesPage.Shapes.Item(5).TextFrame.TextRange.Text = $newDate
    }
}
